# DailyWorks.xlsx
# Bug found while using multiprocessing -> log a "Test 0" entry at 3:18 PM
# as a new row (row 20) at the bottom of the daily schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row -------------------------------------------------
# Column A: the time of the entry (15:18 / 3:18 PM), stored as an Excel
# serial time fraction of a day, same as every other row in column A.
# Column B: the description of the task.
$ws.Cells.Item(20, 1).Value = 0.63750000000000007
$ws.Cells.Item(20, 2).Value = "Test 0"

# --- Match formatting of the existing rows ----------------------------
# Row 18 (and several others) share the same look: a centered time
# value formatted as h:mm AM/PM in column A, and a plain centered label
# in column B. Copy that formatting onto the new row so it fits right in.
$ws.Cells.Item(20, 1).NumberFormat = $ws.Cells.Item(18, 1).NumberFormat
$ws.Cells.Item(20, 1).HorizontalAlignment = $ws.Cells.Item(18, 1).HorizontalAlignment
$ws.Cells.Item(20, 1).VerticalAlignment = $ws.Cells.Item(18, 1).VerticalAlignment

$ws.Cells.Item(20, 2).HorizontalAlignment = $ws.Cells.Item(18, 2).HorizontalAlignment
$ws.Cells.Item(20, 2).VerticalAlignment = $ws.Cells.Item(18, 2).VerticalAlignment

# --- Move the view/selection down to the newly added row --------------
try {
    $excel.ActiveWindow.ScrollRow = 13
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Older/limited hosts may not expose window scrolling; ignore.
}
$ws.Range("B21").Select()
